$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3692786862235664
$ws.Range("C2").Value = 0.06750128626747198
$ws.Range("D2").Value = 0.02643641156462451
$ws.Range("E2").Value = 0.4159123285805748
$ws.Range("F2").Value = 0.709949118744241
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("K2").Value = 0.3603870442761377
$ws.Range("O2").Value = 2.425853481379022
$ws.Range("B3").Value = 0.3267623426023647
$ws.Range("C3").Value = 0.06179291441247869
$ws.Range("D3").Value = 0.02481443763230118
$ws.Range("E3").Value = 0.3628992592812068
$ws.Range("F3").Value = 0.7073304550515829
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("K3").Value = 0.3148674128965183
$ws.Range("O3").Value = 2.431571157312817
$ws.Range("B4").Value = 0.3006692369606014
$ws.Range("C4").Value = 0.05826342446459876
$ws.Range("D4").Value = 0.02381009723379179
$ws.Range("E4").Value = 0.3304359795849052
$ws.Range("F4").Value = 0.7062327119440397
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("K4").Value = 0.2868644905347253
$ws.Range("O4").Value = 2.436779632266507
$ws.Range("B5").Value = 0.2900395173022901
$ws.Range("C5").Value = 0.05681901421840507
$ws.Range("D5").Value = 0.02339872983098701
$ws.Range("E5").Value = 0.3172271342782267
$ws.Range("F5").Value = 0.7059133797577957
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("K5").Value = 0.2754399761406887
$ws.Range("O5").Value = 2.43932806415944
$ws.Range("B6").Value = 0.2882746793737851
$ws.Range("C6").Value = 0.05657880328006115
$ws.Range("D6").Value = 0.0233302972788465
$ws.Range("E6").Value = 0.3150349875119076
$ws.Range("F6").Value = 0.7058680783579092
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("K6").Value = 0.2735421652579078
$ws.Range("O6").Value = 2.439776925921734
$ws.Range("B7").Value = 0.3005258664248061
$ws.Range("C7").Value = 0.05824396932865739
$ws.Range("D7").Value = 0.02380455781108282
$ws.Range("E7").Value = 0.3302577608487525
$ws.Range("F7").Value = 0.7062278873820702
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("K7").Value = 0.2867104680245234
$ws.Range("O7").Value = 2.436812278036115
$ws.Range("B8").Value = 0.354616758735574
$ws.Range("C8").Value = 0.06553816001884627
$ws.Range("D8").Value = 0.02587892582811691
$ws.Range("E8").Value = 0.3976143566367512
$ws.Range("F8").Value = 0.7089401619904834
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("K8").Value = 0.3447032443233979
$ws.Range("O8").Value = 2.427471956428946
$ws.Range("B9").Value = 0.4607732036890297
$ws.Range("C9").Value = 0.07964579805047833
$ws.Range("D9").Value = 0.02987854830372783
$ws.Range("E9").Value = 0.5304715718472153
$ws.Range("F9").Value = 0.7183200586844123
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("K9").Value = 0.4579889468427041
$ws.Range("O9").Value = 2.422674350482225
$ws.Range("B10").Value = 0.5388102127713523
$ws.Range("C10").Value = 0.08988985053686349
$ws.Range("D10").Value = 0.03277412581027761
$ws.Range("E10").Value = 0.6286731617452972
$ws.Range("F10").Value = 0.7277079534868562
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("K10").Value = 0.5409453036401999
$ws.Range("O10").Value = 2.42746075108019
$ws.Range("B11").Value = 0.5743197598488052
$ws.Range("C11").Value = 0.09452372209362636
$ws.Range("D11").Value = 0.03408180682868789
$ws.Range("E11").Value = 0.6735028481930954
$ws.Range("F11").Value = 0.7325254296818002
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("K11").Value = 0.5786236393822719
$ws.Range("O11").Value = 2.431458261238816
$ws.Range("B12").Value = 0.5877674948822857
$ws.Range("C12").Value = 0.09627464524049856
$ws.Range("D12").Value = 0.03457559308593261
$ws.Range("E12").Value = 0.6905033027698551
$ws.Range("F12").Value = 0.7344286564716498
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("K12").Value = 0.5928827030608659
$ws.Range("O12").Value = 2.433234931645643
$ws.Range("B13").Value = 0.5848712448264166
$ws.Range("C13").Value = 0.09589772354327408
$ws.Range("D13").Value = 0.03446931036066303
$ws.Range("E13").Value = 0.6868408354261533
$ws.Range("F13").Value = 0.7340152460770497
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("K13").Value = 0.5898121613913361
$ws.Range("O13").Value = 2.432840582747644
$ws.Range("B14").Value = 0.5754260946170007
$ws.Range("C14").Value = 0.09466784853978538
$ws.Range("D14").Value = 0.03412245925356672
$ws.Range("E14").Value = 0.6749009858912558
$ws.Range("F14").Value = 0.7326804251625276
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("K14").Value = 0.579796922035996
$ws.Range("O14").Value = 2.431599153391147
$ws.Range("B15").Value = 0.5696407906201557
$ws.Range("C15").Value = 0.09391401501758878
$ws.Range("D15").Value = 0.03390981908312796
$ws.Range("E15").Value = 0.6675907181106595
$ws.Range("F15").Value = 0.7318731001349192
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("K15").Value = 0.5736611264519524
$ws.Range("O15").Value = 2.430873014034717
$ws.Range("B16").Value = 0.5364897484143967
$ws.Range("C16").Value = 0.08958648412581738
$ws.Range("D16").Value = 0.03268847110454232
$ws.Range("E16").Value = 0.6257467554810603
$ws.Range("F16").Value = 0.7274041513733636
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("K16").Value = 0.5384817130795625
$ws.Range("O16").Value = 2.427236229955895
$ws.Range("B17").Value = 0.5161549529779563
$ws.Range("C17").Value = 0.0869249336540463
$ws.Range("D17").Value = 0.03193674841050864
$ws.Range("E17").Value = 0.6001184060109495
$ws.Range("F17").Value = 0.7248028866342366
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("K17").Value = 0.5168848815783917
$ws.Range("O17").Value = 2.425472192731064
$ws.Range("B18").Value = 0.5044598935531042
$ws.Range("C18").Value = 0.08539161750438495
$ws.Range("D18").Value = 0.03150348220756172
$ws.Range("E18").Value = 0.5853923531771414
$ws.Range("F18").Value = 0.7233581637361084
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("K18").Value = 0.5044574342174997
$ws.Range("O18").Value = 2.424628816657815
$ws.Range("B19").Value = 0.500500327919184
$ws.Range("C19").Value = 0.08487204175072804
$ws.Range("D19").Value = 0.03135663299388369
$ws.Range("E19").Value = 0.5804088339548485
$ws.Range("F19").Value = 0.7228778334412453
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("K19").Value = 0.5002487790344503
$ws.Range("O19").Value = 2.424372641221396
$ws.Range("B20").Value = 0.5183195285642057
$ws.Range("C20").Value = 0.0872085155361475
$ws.Range("D20").Value = 0.0320168634676179
$ws.Range("E20").Value = 0.6028450535715848
$ws.Range("F20").Value = 0.7250744679569721
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("K20").Value = 0.5191844757123363
$ws.Range("O20").Value = 2.425642245181791
$ws.Range("B21").Value = 0.5782003392402828
$ws.Range("C21").Value = 0.09502919684963729
$ws.Range("D21").Value = 0.03422437620510976
$ws.Range("E21").Value = 0.6784073317457597
$ws.Range("F21").Value = 0.7330703491592914
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("K21").Value = 0.5827388856738764
$ws.Range("O21").Value = 2.431956646954205
$ws.Range("B22").Value = 0.6173417845361655
$ws.Range("C22").Value = 0.1001181529369148
$ws.Range("D22").Value = 0.03565891356750939
$ws.Range("E22").Value = 0.7279348686650167
$ws.Range("F22").Value = 0.7387564645732709
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("K22").Value = 0.6242232693229539
$ws.Range("O22").Value = 2.437616399496335
$ws.Range("B23").Value = 0.5964508569768157
$ws.Range("C23").Value = 0.09740414103174544
$ws.Range("D23").Value = 0.03489403523903434
$ws.Range("E23").Value = 0.7014873838633804
$ws.Range("F23").Value = 0.7356794526076555
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("K23").Value = 0.6020871807872368
$ws.Range("O23").Value = 2.434455034505817
$ws.Range("B24").Value = 0.5173409373524294
$ws.Range("C24").Value = 0.08708031799307037
$ws.Range("D24").Value = 0.03198064684463731
$ws.Range("E24").Value = 0.6016123113793128
$ws.Range("F24").Value = 0.724951527886347
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("K24").Value = 0.5181448638791721
$ws.Range("O24").Value = 2.425564832524003
$ws.Range("B25").Value = 0.43204684069525
$ws.Range("C25").Value = 0.07585044468443414
$ws.Range("D25").Value = 0.02880398899296921
$ws.Range("E25").Value = 0.4944351383057182
$ws.Range("F25").Value = 0.715345636987621
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("K25").Value = 0.4273898620322711
$ws.Range("O25").Value = 2.422517752854048
